# Apply the updated cryptocurrency price/volume figures for the
# "Updated cryptos list" GitHub Actions refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text looks like a plain decimal number (e.g. "1.000",
# "21.29") are forced to Text format first so Excel does not silently
# reinterpret them as numeric values (they are inline text in the sheet).
$numberLikeUpdates = @{
    "D5" = "331.15"
    "D6" = "1.000"
    "D7" = "0.4624"
    "D9" = "47.97"
    "D10" = "0.07872"
    "D11" = "0.9865"
    "D12" = "21.29"
    "D14" = "5.853"
    "D15" = "7.007"
    "D17" = "88.24"
    "D18" = "0.06533"
    "D23" = "5.354"
    "D24" = "10.92"
    "D25" = "2.233"
    "D27" = "156.93"
    "D28" = "19.36"
    "D29" = "2.064"
    "D30" = "5.312"
    "D31" = "117.68"
    "D32" = "0.9608"
    "D33" = "0.09358"
    "D36" = "5.257"
    "D37" = "0.06038"
    "D38" = "0.02203"
    "D39" = "8.286"
    "D40" = "1.163"
    "D41" = "1.000"
    "D42" = "0.5752"
    "D44" = "10.05"
    "D45" = "1.267"
    "D46" = "2.298"
    "D47" = "11.95"
    "D48" = "0.5436"
    "D49" = "0.07159"
    "D50" = "1.893"
    "D51" = "111.56"
}

# Remaining cells (multi-dot price strings and the "  +x.xx%  " volume
# strings) are never parsed as numbers, so they can be set directly.
$textUpdates = @{
    "D2" = "28.407.20"
    "E2" = "  +0.22%  "
    "D3" = "1.867.34"
    "E3" = "  -0.10%  "
    "E4" = "  -0.03%  "
    "E5" = "  -2.25%  "
    "E6" = "  -0.05%  "
    "E7" = "  -1.70%  "
    "E8" = "  +1.53%  "
    "E9" = "  +1.22%  "
    "E10" = "  -1.50%  "
    "E11" = "  -2.11%  "
    "E12" = "  -3.02%  "
    "D13" = "1.862.88"
    "E13" = "  -0.67%  "
    "E14" = "  -2.60%  "
    "E15" = "  -3.78%  "
    "E16" = "  -0.09%  "
    "E17" = "  -3.22%  "
    "E18" = "  -1.09%  "
    "E19" = "  -1.75%  "
    "E20" = "  -2.73%  "
    "E21" = "  +0.04%  "
    "D22" = "28.389.16"
    "E22" = "  +0.09%  "
    "E23" = "  -2.00%  "
    "E24" = "  -1.35%  "
    "E25" = "  -2.60%  "
    "D26" = "2.083.84"
    "E26" = "  -0.66%  "
    "E27" = "  -1.95%  "
    "E28" = "  -2.58%  "
    "E29" = "  -4.09%  "
    "E30" = "  -3.31%  "
    "E31" = "  -2.19%  "
    "E32" = "  -1.47%  "
    "E33" = "  -1.66%  "
    "E34" = "  -0.20%  "
    "E35" = "  +0.13%  "
    "E36" = "  -1.90%  "
    "E37" = "  -0.99%  "
    "E38" = "  -3.15%  "
    "E39" = "  -2.24%  "
    "E40" = "  -1.61%  "
    "E41" = "  -0.06%  "
    "E42" = "  -3.89%  "
    "E43" = "  -4.00%  "
    "E44" = "  -3.34%  "
    "E45" = "  -2.87%  "
    "E46" = "  +13.52%  "
    "E47" = "  -1.67%  "
    "E48" = "  -3.38%  "
    "E49" = "  +3.37%  "
    "E50" = "  -3.86%  "
    "E51" = "  +0.24%  "
}

foreach ($cell in $numberLikeUpdates.Keys) {
    $range = $ws.Range($cell)
    $range.NumberFormat = "@"
    $range.Value = $numberLikeUpdates[$cell]
    $range.ClearFormats()
}

foreach ($cell in $textUpdates.Keys) {
    $ws.Range($cell).Value = $textUpdates[$cell]
}
